# Optimizer.pptx — "Fix English on Goals slide; constraint management."
#
# Uses TextRange.Characters(start, length) sub-ranges (rather than whole
# paragraph .Text assignment) so that only the targeted characters are
# rewritten and PowerPoint's word-level diffing doesn't re-split runs
# that should stay untouched / doesn't keep the wrong run's rPr when two
# runs get merged into one.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 — "Goals of this session" bullet: reword.
# ---------------------------------------------------------------------
$s3   = $p.Slides.Item(3)
$tr3  = $s3.Shapes.Item(2).TextFrame.TextRange
$para3 = $tr3.Paragraphs(7, 1)
$old3 = "Understanding to help write better queries"
$i3   = $para3.Text.IndexOf($old3)
$para3.Characters($i3 + 1, $old3.Length).Text = "Deeper understanding: write better queries!"

# ---------------------------------------------------------------------
# Slide 35 — "Extensive use of heuristics to prune search space": the
# trailing "search spac" + "e" runs become a single "search space" run.
# ---------------------------------------------------------------------
$s35   = $p.Slides.Item(35)
$tr35  = $s35.Shapes.Item(2).TextFrame.TextRange
$para35 = $tr35.Paragraphs(8, 1)
$old35 = "search spac"
$i35   = $para35.Text.IndexOf($old35)
# Rewrite the "search spac" run to the full merged word, keeping its rPr.
$para35.Characters($i35 + 1, $old35.Length).Text = "search space"
# The former trailing "e" run now sits right after the new text; empty it
# out so it doesn't leave a duplicated "search spacee".
$para35b = $tr35.Paragraphs(8, 1)
$newEnd35 = $i35 + ("search space").Length
$para35b.Characters($newEnd35 + 1, 1).Text = ""

# ---------------------------------------------------------------------
# Slide 55 — "The Optimizer Is Exceptionally Complex" content body.
# ---------------------------------------------------------------------
$s55  = $p.Slides.Item(55)
$tr55 = $s55.Shapes.Item(2).TextFrame.TextRange

# Halloween protection bullet: append ", constraint management".
$para55a = $tr55.Paragraphs(3, 1)
$old55a  = "Halloween protection, triggers, index updates"
$i55a    = $para55a.Text.IndexOf($old55a)
$para55a.Characters($i55a + 1, $old55a.Length).Text = "Halloween protection, triggers, index updates, constraint management"

# Window functions bullet: merge the "Window functions, partitioned " and
# "tables, " runs into a single run, keeping the *second* run's rPr
# (sz="2200" dirty="0") and dropping the first run entirely.
$para55b   = $tr55.Paragraphs(7, 1)
$oldTables = "tables, "
$iTables   = $para55b.Text.IndexOf($oldTables)
$para55b.Characters($iTables + 1, $oldTables.Length).Text = "Window functions, partitioned tables, "
$para55c  = $tr55.Paragraphs(7, 1)
$oldLead  = "Window functions, partitioned "
$iLead    = $para55c.Text.IndexOf($oldLead)
$para55c.Characters($iLead + 1, $oldLead.Length).Text = ""

# ---------------------------------------------------------------------
# Slide 62 — "Thank You" links list: rename the Parse Tree Viewer tool.
# ---------------------------------------------------------------------
$s62  = $p.Slides.Item(62)
$tr62 = $s62.Shapes.Item(2).TextFrame.TextRange
$para62 = $tr62.Paragraphs(5, 1)
$old62  = "SQL Server Parse Tree Viewer binaries & source"
$i62    = $para62.Text.IndexOf($old62)
$para62.Characters($i62 + 1, $old62.Length).Text = "SQL Server Query Tree Viewer binaries & source"
